$wb = $excel.ActiveWorkbook

# --- Metadata sheet (sheet1.xml) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/description"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet (sheet2.xml) ---
$elements = $wb.Worksheets.Item("Elements")
# Fixed Value for Extension.url (row 5) picks up the new canonical URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/description"
# Constraint(s) for the root Extension row (row 2) is cleared
$elements.Range("AI2").Value = ""
